$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns E (5) through T (20), for rows 2 through 17.
# Each inner array corresponds to one row, in order E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T
$data = @(
    @(3,1,40.34291466666667,121.028744,0.3404392602027053,0.3404392602027053,3,1,22.51161266666667,67.53483800000001,0.3150979864474181,0.3150979864474181,908.1840688203858,8173.656619383473,0.1072717253975211,0.1072717253975211),
    @(3,1,40.34291466666667,121.028744,0.3404392602027053,0.3404392602027053,3,1,30.248441,90.745323,0.4233913844114141,0.4233913844114141,1220.310274062701,10982.79246656431,0.144139049685221,0.144139049685221),
    @(3,1,40.34291466666667,121.028744,0.3404392602027053,0.3404392602027053,3,1,18.63333,55.89999,0.2608131567803688,0.2608131567803688,751.72284214584,6765.50557931256,0.08879103814544094,0.08879103814544094),
    @(3,1,40.34291466666667,121.028744,0.3404392602027053,0.3404392602027053,2,0.6666666666666666,0.04982966666666667,0.149489,0.000697472360799001,0.000697472360799001,2.010273990201778,18.092465911816,0.0002374469745222462,0.0002374469745222463),
    @(3,1,41.443863,124.331589,0.3497297648481489,0.3497297648481489,3,1,22.51161266666667,67.53483800000001,0.3150979864474181,0.3150979864474181,932.9681912663981,8396.713721397582,0.1101991447043807,0.1101991447043807),
    @(3,1,41.443863,124.331589,0.3497297648481489,0.3497297648481489,3,1,30.248441,90.745323,0.4233913844114141,0.4233913844114141,1253.612244767583,11282.51020290825,0.1480725693089361,0.1480725693089361),
    @(3,1,41.443863,124.331589,0.3497297648481489,0.3497297648481489,3,1,18.63333,55.89999,0.2608131567803688,0.2608131567803688,772.2371757537901,6950.13458178411,0.09121412399010177,0.09121412399010176),
    @(3,1,41.443863,124.331589,0.3497297648481489,0.3497297648481489,2,0.6666666666666666,0.04982966666666667,0.149489,0.000697472360799001,0.000697472360799001,2.065133878669,18.586204908021,0.0002439268447303179,0.0002439268447303179),
    @(3,1,36.642055,109.926165,0.3092090445020276,0.3092090445020277,3,1,22.51161266666667,67.53483800000001,0.3150979864474181,0.3150979864474181,824.8717494706967,7423.845745236271,0.097431147313919,0.09743114731391903),
    @(3,1,36.642055,109.926165,0.3092090445020276,0.3092090445020277,3,1,30.248441,90.745323,0.4233913844114141,0.4233913844114141,1108.365038786255,9975.285349076295,0.130916445424244,0.130916445424244),
    @(3,1,36.642055,109.926165,0.3092090445020276,0.3092090445020277,3,1,18.63333,55.89999,0.2608131567803688,0.2608131567803688,682.76350269315,6144.871524238351,0.08064578700161537,0.08064578700161537),
    @(3,1,36.642055,109.926165,0.3092090445020276,0.3092090445020277,2,0.6666666666666666,0.04982966666666667,0.149489,0.000697472360799001,0.000697472360799001,1.825861386631667,16.432752479685,0.0002156647622492326,0.0002156647622492326),
    @(3,1,0.07370033333333333,0.221101,0.0006219304471182344,0.0006219304471182345,3,1,22.51161266666667,67.53483800000001,0.3150979864474181,0.3150979864474181,1.659113357404222,14.932020216638,0.0001959690315972981,0.0001959690315972981),
    @(3,1,0.07370033333333333,0.221101,0.0006219304471182344,0.0006219304471182345,3,1,30.248441,90.745323,0.4233913844114141,0.4233913844114141,2.229320184513667,20.063881660623,0.000263319993012999,0.0002633199930129991),
    @(3,1,0.07370033333333333,0.221101,0.0006219304471182344,0.0006219304471182345,3,1,18.63333,55.89999,0.2608131567803688,0.2608131567803688,1.37328263211,12.35954368899,0.000162207643210733,0.0001622076432107329),
    @(3,1,0.07370033333333333,0.221101,0.0006219304471182344,0.0006219304471182345,2,0.6666666666666666,0.04982966666666667,0.149489,0.000697472360799001,0.000697472360799001,0.003672463043222222,0.033052167389,0.0000004337792972043332,0.0000004337792972043333)
)

$startRow = 2
$startCol = 5

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($startRow + $i, $startCol + $j).Value = $rowVals[$j]
    }
}
